$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("GNG_TO-16502912950093892").Name = "GNG_TO-165047789871423"
$wb.Worksheets.Item("NB_TO-1650291297313209").Name = "NB_TO-1650477900315923"
$wb.Worksheets.Item("RS_TO-16502912973142076").Name = "RS_TO-16504779003169143"
$wb.Worksheets.Item("TOL_TO-16502912973773").Name = "TOL_TO-1650477900364918"
$wb.Worksheets.Item("vSAT_TO-16502912974481826").Name = "vSAT_TO-16504779004279475"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item("GNG_TO-165047789871423")
$ws1.Range("B2").Value = "go_stims-16504778986792343.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778986972299.csv"
$ws1.Range("B4").Value = "go_stims-1650477898699232.csv"
$ws1.Range("B5").Value = "GNG_stims-1650477898713265.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item("NB_TO-1650477900315923")
$ws2.Range("B2").Value = "TB-16504779002922618.csv"
$ws2.Range("B3").Value = "OB-16504778988772328.csv"
$ws2.Range("B4").Value = "OB-16504778999622335.csv"
$ws2.Range("B5").Value = "TB-16504779001872308.csv"
$ws2.Range("B6").Value = "TB-16504779000912633.csv"
$ws2.Range("B7").Value = "ZB-match_5-1650477898743265.csv"
$ws2.Range("B8").Value = "ZB-match_1-16504778988442616.csv"
$ws2.Range("B9").Value = "OB-16504778990132678.csv"
$ws2.Range("B10").Value = "ZB-match_1-16504778987792604.csv"

# Sheet 3 (RS)
$ws3 = $wb.Worksheets.Item("RS_TO-16504779003169143")
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item("TOL_TO-1650477900364918")
$ws4.Range("B2").Value = "MM_stims-1650477900332042.csv"
$ws4.Range("B3").Value = "ZM_stims-16504779003199496.csv"
$ws4.Range("B4").Value = "MM_stims-1650477900347952.csv"
$ws4.Range("B5").Value = "ZM_stims-16504779003329177.csv"
$ws4.Range("B6").Value = "MM_stims-1650477900363914.csv"
$ws4.Range("B7").Value = "ZM_stims-1650477900347952.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item("vSAT_TO-16504779004279475")
$ws5.Range("B2").Value = "SAT_stims-16504779003669193.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504779003959498.csv"
$ws5.Range("B4").Value = "SAT_stims-16504779003809636.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504779004119487.csv"
